$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Specification")

# For rows 2 through 20, the values previously held in columns C:E need to
# shift one column to the right (C->D, D->E, E->F), and column C gets the
# new value "Application". Column F was always empty in these rows before
# the edit (the real content stopped at E at most), and column G already
# holds the "description" text which must be left untouched. Process the
# shift from the rightmost column down to the leftmost so we don't clobber
# a value before it is read.

for ($r = 2; $r -le 20; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2

    $ws.Cells.Item($r, 6).Value = $eVal
    $ws.Cells.Item($r, 5).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $cVal
    $ws.Cells.Item($r, 3).Value = "Application"
}
